$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.130.97"
$ws.Range("E2").Value = "  +0.21%  "

$ws.Range("D3").Value = "1.841.14"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("D4").Value = "'0.9987"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'243.79"
$ws.Range("E5").Value = "  -0.22%  "

$ws.Range("D6").Value = "'0.6249"
$ws.Range("E6").Value = "  -1.35%  "

$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "'0.07527"
$ws.Range("E8").Value = "  -0.75%  "

$ws.Range("E10").Value = "  +2.53%  "

$ws.Range("D11").Value = "'0.07707"
$ws.Range("E11").Value = "  -0.51%  "

$ws.Range("D12").Value = "1.838.49"
$ws.Range("E12").Value = "  +0.17%  "

$ws.Range("D13").Value = "'5.029"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("D14").Value = "'0.6773"
$ws.Range("E14").Value = "  +0.83%  "

$ws.Range("D15").Value = "'83.15"
$ws.Range("E15").Value = "  -0.07%  "

$ws.Range("D16").Value = "'0.000009340"
$ws.Range("E16").Value = "  -5.29%  "

$ws.Range("D17").Value = "'5.982"
$ws.Range("E17").Value = "  -2.31%  "

$ws.Range("D18").Value = "29.128.32"
$ws.Range("E18").Value = "  +0.12%  "

$ws.Range("D19").Value = "2.083.84"
$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Value = "'12.72"
$ws.Range("E20").Value = "  +1.24%  "

$ws.Range("D21").Value = "'229.19"
$ws.Range("E21").Value = "  +0.94%  "

$ws.Range("D22").Value = "'1.001"
$ws.Range("E22").Value = "  +0.05%  "

$ws.Range("D23").Value = "'7.175"
$ws.Range("E23").Value = "  -0.76%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("E24").Value = "  +0.08%  "

$ws.Range("D25").Value = "'160.45"
$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").Value = "'0.1402"
$ws.Range("E26").Value = "  -0.25%  "

$ws.Range("D27").Value = "'8.555"
$ws.Range("E27").Value = "  +0.05%  "

$ws.Range("D28").Value = "'17.95"
$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("D29").Value = "'1.498"
$ws.Range("E29").Value = "  -0.35%  "

$ws.Range("D30").Value = "'4.191"
$ws.Range("E30").Value = "  +1.64%  "

$ws.Range("D31").Value = "'4.149"
$ws.Range("E31").Value = "  +2.20%  "

$ws.Range("D32").Value = "'0.05567"
$ws.Range("E32").Value = "  +3.14%  "

$ws.Range("E33").Value = "  +0.22%  "

$ws.Range("D34").Value = "'0.7497"
$ws.Range("E34").Value = "  +0.32%  "

$ws.Range("D35").Value = "'1.851"
$ws.Range("E35").Value = "  -0.66%  "

$ws.Range("D36").Value = "'1.149"
$ws.Range("E36").Value = "  +0.52%  "

$ws.Range("D37").Value = "'2.667"
$ws.Range("E37").Value = "  -0.17%  "

$ws.Range("D38").Value = "1.236.25"
$ws.Range("E38").Value = "  -0.74%  "

$ws.Range("D39").Value = "'2.775"
$ws.Range("E39").Value = "  +0.49%  "

$ws.Range("D40").Value = "'0.01788"
$ws.Range("E40").Value = "  -0.46%  "

$ws.Range("D41").Value = "'6.605"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("D42").Value = "'0.9023"
$ws.Range("E42").Value = "  -0.46%  "

$ws.Range("D43").Value = "'1.000"
$ws.Range("E43").Value = "  -0.18%  "

$ws.Range("D44").Value = "'102.53"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").Value = "1.983.41"
$ws.Range("E45").Value = "  -0.33%  "

$ws.Range("D46").Value = "'66.56"
$ws.Range("E46").Value = "  +2.66%  "

$ws.Range("D47").Value = "'0.00000000122"
$ws.Range("E47").Value = "  -0.87%  "

$ws.Range("D48").Value = "'0.5088"
$ws.Range("E48").Value = "  -0.59%  "

$ws.Range("D49").Value = "'0.4088"
$ws.Range("E49").Value = "  -0.47%  "

$ws.Range("D50").Value = "'9.093"
$ws.Range("E50").Value = "  +0.83%  "

$ws.Range("B51").Value = "XinFinNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D51").Value = "'0.07266"
$ws.Range("E51").Value = "  +17.09%  "
